# QA fixes: normalize "Creation date" timestamps from 4-digit to 2-digit
# years, and append three newly-coded segments (Location:City,
# Location:Country, Location:Hospital name) for document 22121.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Shorten the 4-digit year in every "Creation date" (column M) cell ---
$dates = $ws.Range("M1:M183")
$dates.Replace("/2018 ", "/18 ")
$dates.Replace("/2019 ", "/19 ")

# --- 2. Append three new coded-segment rows (184-186) ---
# Copy the formatting (styles/fills/borders) of the last existing data row
# down onto the new rows first, then fill in the values.
$ws.Range("A183:M183").Copy()
$ws.Range("A184:M186").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$rows = @(
    @{ Row=184; Code="22121"; CodeName="Location:City";          Begin="1: 1314"; End="1: 1319"; Segment="Taipei"; Area=6;  Coverage=0.041531113726033089; Author="emmamendelsohn"; Created="8/22/19 14:19:16" },
    @{ Row=185; Code="22121"; CodeName="Location:Country";       Begin="1: 1326"; End="1: 1331"; Segment="Taiwan"; Area=6;  Coverage=0.041531113726033089; Author="emmamendelsohn"; Created="8/22/19 14:19:20" },
    @{ Row=186; Code="22121"; CodeName="Location:Hospital name"; Begin="1: 1213"; End="1: 1275"; Segment="Tri-Service General Hospital,  `nNational Defense Medical Center"; Area=62; Coverage=0.4291548418356752; Author="emmamendelsohn"; Created="8/22/19 14:19:48" }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = "●"
    $ws.Cells.Item($row, 2).Value = ""
    $ws.Cells.Item($row, 3).Value = ""
    $ws.Cells.Item($row, 4).Value = $r.Code
    $ws.Cells.Item($row, 5).Value = $r.CodeName
    $ws.Cells.Item($row, 6).Value = $r.Begin
    $ws.Cells.Item($row, 7).Value = $r.End
    $ws.Cells.Item($row, 8).Value = 0
    $ws.Cells.Item($row, 9).Value = $r.Segment
    $ws.Cells.Item($row, 10).Value = $r.Area
    $ws.Cells.Item($row, 11).Value = $r.Coverage
    $ws.Cells.Item($row, 12).Value = $r.Author
    $ws.Cells.Item($row, 13).Value = $r.Created
}

# --- 3. Row heights for the new rows (single-line vs. wrapped text) ---
$ws.Rows.Item(184).RowHeight = 16
$ws.Rows.Item(185).RowHeight = 16
$ws.Rows.Item(186).RowHeight = 30

Write-Output "done"
